# Capabilities mobile section 1
#
# Adds a new introductory row (row 3) to the CAPABILITIES sheet with a
# "SECTION_1_TEXT_1_SMALL"-style label in column A and a rich-text
# paragraph in column B, then leaves CAPABILITIES as the active sheet
# with B4 selected (ABOUT, previously active, loses its tab selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAPABILITIES")

# Insert a new blank row above the current row 4 (pushes everything else
# down by one and inherits formatting from the row above, like Excel does).
$ws.Rows("3:3").Insert()

# Column A: plain label string.
$ws.Range("A3").Value = "SECTION_1_TEXT_1_SMALL"

# Column B: rich text with a couple of inline-colored "runs" mimicking the
# <span className='{{custom}}'>...</span> convention used throughout the
# workbook's copy.
$text = "<span className='{{custom}}'>Good design</span>`nmake perfect sense`nto your goal."
$ws.Range("B3").Value = $text

# "custom"                -> green
$ws.Range("B3").Characters(20, 6).Font.Color = 4697456
# "}}'>"                  -> blue
$ws.Range("B3").Characters(26, 4).Font.Color = 12874308
# "Good design"           -> orange
$ws.Range("B3").Characters(30, 11).Font.Color = 3243501
# "</span>" + newline     -> blue
$ws.Range("B3").Characters(41, 8).Font.Color = 12874308

# Match the row height used by the other "big" rows in this sheet.
$ws.Rows("3:3").RowHeight = 45

# Make CAPABILITIES the active sheet/tab, with B4 as the selected cell
# (mirrors having just typed the text above and pressed Enter).
$ws.Activate() | Out-Null
$ws.Range("B4").Select() | Out-Null

Write-Output "Inserted CAPABILITIES row 3 and activated CAPABILITIES sheet"
